# Rename settings to use underscores, and simplify the countries_csv
# choice_filter formula on the "survey" sheet.

$wb = $excel.ActiveWorkbook

# --- survey sheet: simplify the choice_filter for select_one countries_csv ---
$survey = $wb.Worksheets.Item("survey")
$survey.Range("I18").Value = "context.region === data('region')"

# --- settings sheet: rename settings to use underscores ---
$settings = $wb.Worksheets.Item("settings")
$settings.Range("A2").Value = "form_id"
$settings.Range("A3").Value = "form_version"
$settings.Range("A4").Value = "form_title"
